$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("December", 18, "00:00:05", 1.88, 0, 45, 1032, 2.26)
    ,@("December", 18, "00:30:05", 1.88, 0, 47, 1031, 2.5)
    ,@("December", 18, "01:00:05", 1.81, 0, 84, 1031, 2.48)
    ,@("December", 18, "01:30:05", 1.63, 0, 83, 1031, 2.42)
    ,@("December", 18, "02:00:04", 1.38, 0, 51, 1031, 2.42)
    ,@("December", 18, "02:30:05", 1.31, 0, 54, 1031, 1.89)
    ,@("December", 18, "03:00:04", 1.19, 0, 85, 1031, 1.89)
    ,@("December", 18, "03:30:05", 1.31, 0, 86, 1031, 2.68)
    ,@("December", 18, "04:00:06", 1.5, 0, 81, 1031, 2.24)
    ,@("December", 18, "04:30:05", 1.63, 0, 81, 1031, 2.24)
    ,@("December", 18, "05:00:04", 1.63, 0, 80, 1031, 1.34)
    ,@("December", 18, "05:30:05", 2.19, 0, 85, 1031, 3.13)
    ,@("December", 18, "06:00:04", 3.19, 0, 83, 1031, 2.24)
    ,@("December", 18, "06:30:05", 3.81, 0, 82, 1031, 1.34)
    ,@("December", 18, "07:00:05", 3.88, 0, 80, 1031, 2.68)
    ,@("December", 18, "07:30:04", 4.31, 0, 82, 1031, 2.24)
    ,@("December", 18, "08:00:04", 4.5, 0, 80, 1031, 0.89)
    ,@("December", 18, "08:30:05", 4.38, 0, 81, 1031, 1.79)
    ,@("December", 18, "09:00:04", 4.19, 0, 79, 1031, 1.79)
    ,@("December", 18, "09:30:05", 4.38, 0, 80, 1031, 3.13)
    ,@("December", 18, "10:00:04", 4.88, 0, 79, 1031, 2.68)
    ,@("December", 18, "10:30:05", 5.69, 0, 77, 1031, 3.58)
    ,@("December", 18, "11:00:05", 6.69, 0, 69, 1031, 0.89)
    ,@("December", 18, "11:30:05", 7.81, 0, 74, 1030, 2.68)
    ,@("December", 18, "12:00:05", 9, 0, 63, 1030, 1.79)
    ,@("December", 18, "12:30:04", 10, 0, 70, 1028, 1.79)
    ,@("December", 18, "13:00:05", 10.88, 0, 69, 1029, 2.24)
    ,@("December", 18, "13:30:04", 11.63, 0, 67, 1027, 2.68)
    ,@("December", 18, "14:00:05", 12.19, 0, 65, 1027, 3.58)
    ,@("December", 18, "14:30:04", 12.69, 0, 65, 1026, 3.58)
    ,@("December", 18, "15:00:05", 13.13, 0, 65, 1026, 3.13)
    ,@("December", 18, "15:30:04", 13.31, 0, 66, 1025, 3.58)
    ,@("December", 18, "16:00:05", 13.19, 0, 68, 1025, 3.13)
    ,@("December", 18, "16:30:05", 12.88, 0, 68, 1025, 2.68)
    ,@("December", 18, "17:00:05", 12.38, 0, 68, 1025, 1.34)
    ,@("December", 18, "17:30:05", 12.13, 0, 70, 1025, 0.45)
    ,@("December", 18, "18:00:05", 10.88, 0, 72, 1025, 2.68)
    ,@("December", 18, "18:30:05", 10.38, 0, 73, 1025, 1.79)
    ,@("December", 18, "19:00:05", 10.13, 0, 77, 1025, 3.13)
    ,@("December", 18, "19:30:04", 9.19, 0, 75, 1025, 4.02)
)

$startRow = 82
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

Write-Host $ws.UsedRange.Address()
